$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 18 de Abril de 2020 a las 14:52"

# Alemania (row 8) - refreshed case numbers, no reordering
$ws.Range("B8").Value = 141483
$ws.Range("C8").Value = 86
$ws.Range("E8").Value = 51721
$ws.Range("G8").Value = 10
$ws.Range("H8").Value = 4362

# India overtakes Austria (row 20/21 swap with India getting new data)
$ws.Range("A20").Value = "India"
$ws.Range("B20").Value = 14792
$ws.Range("C20").Value = 440
$ws.Range("D20").Value = 2045
$ws.Range("E20").Value = 12259
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 488

$ws.Range("A21").Value = "Austria"
$ws.Range("B21").Value = 14671
$ws.Range("C21").Value = 76
$ws.Range("D21").Value = 10214
$ws.Range("E21").Value = 4014
$ws.Range("F21").Value = 208
$ws.Range("G21").Value = 12
$ws.Range("H21").Value = 443

# Arabia Saudita overtakes Pakistan and Dinamarca (rows 32/33/34 shift)
$ws.Range("A32").Value = "Arabia Saudita"
$ws.Range("B32").Value = 8274
$ws.Range("C32").Value = 1132
$ws.Range("D32").Value = 1329
$ws.Range("E32").Value = 6853
$ws.Range("F32").Value = 74
$ws.Range("G32").Value = 5
$ws.Range("H32").Value = 92

$ws.Range("A33").Value = "Pakistan"
$ws.Range("B33").Value = 7638
$ws.Range("C33").Value = 613
$ws.Range("D33").Value = 1832
$ws.Range("E33").Value = 5663
$ws.Range("F33").Value = 46
$ws.Range("G33").Value = 8
$ws.Range("H33").Value = 143

$ws.Range("A34").Value = "Dinamarca"
$ws.Range("B34").Value = 7242
$ws.Range("C34").Value = 169
$ws.Range("D34").Value = 3847
$ws.Range("E34").Value = 3049
$ws.Range("F34").Value = 76
$ws.Range("G34").Value = 10
$ws.Range("H34").Value = 346

# Kazajistan (row 67) - refreshed case numbers, no reordering
$ws.Range("B67").Value = 1615
$ws.Range("C67").Value = 69
$ws.Range("D67").Value = 377
$ws.Range("E67").Value = 1221
